$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 251.3125
$ws.Range("I39").Value = 250.85715
$ws.Range("J39").Value = 254.5
$ws.Range("K39").Value = 752.5714499999999
$ws.Range("L39").Value = 763.5
$ws.Range("M39").Value = -456.5714499999999
$ws.Range("N39").Value = -1355.5
# Row 43
$ws.Range("H43").Value = 8343557.5
$ws.Range("J43").Value = 35000
$ws.Range("L43").Value = 35000
$ws.Range("N43").Value = -35138
# Row 53
$ws.Range("H53").Value = 453.16666
$ws.Range("I53").Value = 404.75
$ws.Range("K53").Value = 404.75
$ws.Range("M53").Value = 232.25
# Row 92
$ws.Range("H92").Value = 1288.9048
$ws.Range("I92").Value = 1323.1177
$ws.Range("J92").Value = 1143.5
$ws.Range("K92").Value = 1323.1177
$ws.Range("L92").Value = 1143.5
$ws.Range("M92").Value = -75.11770000000001
$ws.Range("N92").Value = -3639.5
# Row 100
$ws.Range("H100").Value = 3511
$ws.Range("I100").Value = 1913
$ws.Range("J100").Value = 7506
$ws.Range("K100").Value = 1913
$ws.Range("L100").Value = 7506
$ws.Range("M100").Value = -1372
$ws.Range("N100").Value = -8588
# Row 116
$ws.Range("H116").Value = 6856.857
$ws.Range("J116").Value = 9000.5
$ws.Range("L116").Value = 9000.5
$ws.Range("N116").Value = -15884.5
# Row 138
$ws.Range("H138").Value = 4195.5
$ws.Range("I138").Value = 3513
$ws.Range("J138").Value = 4726.3335
$ws.Range("K138").Value = 10539
$ws.Range("L138").Value = 14179.0005
$ws.Range("M138").Value = -5399
$ws.Range("N138").Value = -24459.0005

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6349.9
$ws.Range("I61").Value = 4999.8
$ws.Range("K61").Value = 4999.8
$ws.Range("M61").Value = -4787.8
# Row 63
$ws.Range("H63").Value = 5120.6
$ws.Range("I63").Value = 868.3333
$ws.Range("J63").Value = 11499
$ws.Range("K63").Value = 868.3333
$ws.Range("L63").Value = 11499
$ws.Range("M63").Value = -182.3333
$ws.Range("N63").Value = -12871
# Row 66
$ws.Range("H66").Value = 5120.6
$ws.Range("I66").Value = 868.3333
$ws.Range("J66").Value = 11499
$ws.Range("K66").Value = 4341.6665
$ws.Range("L66").Value = 57495
$ws.Range("M66").Value = -909.6665000000003
$ws.Range("N66").Value = -64359
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 122
$ws.Range("H122").Value = 3854.3845
$ws.Range("I122").Value = 3854.3845
$ws.Range("K122").Value = 11563.1535
$ws.Range("M122").Value = -9113.1535
# Row 132
$ws.Range("H132").Value = 2883.1667
$ws.Range("I132").Value = 2883.1667
$ws.Range("K132").Value = 8649.500100000001
$ws.Range("M132").Value = -6119.500100000001
# Row 136
$ws.Range("H136").Value = 6349.9
$ws.Range("I136").Value = 4999.8
$ws.Range("K136").Value = 14999.4
$ws.Range("M136").Value = -12449.4

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 8932
$ws.Range("I82").Value = 8932
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8932
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -8549
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 8932
$ws.Range("I85").Value = 8932
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8932
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -7606
$ws.Range("N85").ClearContents()
# Row 107
$ws.Range("H107").Value = 7022.1665
$ws.Range("I107").Value = 6785.4287
$ws.Range("J107").Value = 7353.6
$ws.Range("K107").Value = 6785.4287
$ws.Range("L107").Value = 7353.6
$ws.Range("M107").Value = -4865.4287
$ws.Range("N107").Value = -11193.6
# Row 134
$ws.Range("H134").Value = 3801.818
$ws.Range("I134").Value = 3801.818
$ws.Range("K134").Value = 11405.454
$ws.Range("M134").Value = -8870.454000000002

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 180
$ws.Range("I8").Value = 500
$ws.Range("J8").Value = 116
$ws.Range("K8").Value = 500
$ws.Range("L8").Value = 116
$ws.Range("M8").Value = -360
$ws.Range("N8").Value = -396
# Row 62
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3248
# Row 65
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 10000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -16240
# Row 122
$ws.Range("H122").Value = 2098.8572
$ws.Range("J122").Value = 994
$ws.Range("L122").Value = 2982
$ws.Range("N122").Value = -7882
# Row 132
$ws.Range("H132").Value = 865.6667
$ws.Range("I132").Value = 865.6667
$ws.Range("K132").Value = 2597.0001
$ws.Range("M132").Value = -67.0001000000002
# Row 134
$ws.Range("H134").Value = 1033.6154
$ws.Range("I134").Value = 1034.96
$ws.Range("K134").Value = 3104.88
$ws.Range("M134").Value = -569.8800000000001

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 197.4375
$ws.Range("I10").Value = 246.81818
$ws.Range("J10").Value = 88.8
$ws.Range("K10").Value = 740.4545400000001
$ws.Range("L10").Value = 266.4
$ws.Range("M10").Value = -601.4545400000001
$ws.Range("N10").Value = -544.4
# Row 14
$ws.Range("H14").Value = 1727.2
$ws.Range("I14").Value = 1727.2
$ws.Range("K14").Value = 5181.6
$ws.Range("M14").Value = -5008.6
# Row 102
$ws.Range("H102").Value = 3000
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 9000
$ws.Range("N102").Value = -13868
# Row 131
$ws.Range("H131").Value = 1705.3214
$ws.Range("I131").Value = 766
$ws.Range("K131").Value = 2298
$ws.Range("M131").Value = 2742

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 6
$ws.Range("H6").Value = 460.33334
$ws.Range("J6").Value = 399.5
$ws.Range("L6").Value = 399.5
$ws.Range("N6").Value = -625.5
# Row 9
$ws.Range("H9").Value = 302.33334
$ws.Range("I9").Value = 203.5
$ws.Range("K9").Value = 203.5
$ws.Range("M9").Value = -33.5
# Row 16
$ws.Range("H16").Value = 460.33334
$ws.Range("J16").Value = 399.5
$ws.Range("L16").Value = 399.5
$ws.Range("N16").Value = -899.5
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# Row 36
$ws.Range("H36").Value = 27039
$ws.Range("I36").Value = 27039
$ws.Range("K36").Value = 27039
$ws.Range("M36").Value = -26554
# Row 70
$ws.Range("H70").Value = 6137.8
$ws.Range("I70").Value = 5894.5
$ws.Range("K70").Value = 5894.5
$ws.Range("M70").Value = -5624.5
# Row 73
$ws.Range("H73").Value = 6137.8
$ws.Range("I73").Value = 5894.5
$ws.Range("K73").Value = 5894.5
$ws.Range("M73").Value = -4958.5
# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
# Row 126
$ws.Range("H126").Value = 4271.8335
$ws.Range("I126").Value = 3803.4
$ws.Range("K126").Value = 11410.2
$ws.Range("M126").Value = -8940.200000000001

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 14646.235
# Row 22
$ws.Range("H22").Value = 1622.2941
$ws.Range("I22").Value = 369.33334
$ws.Range("K22").Value = 369.33334
$ws.Range("M22").Value = -74.33334000000002
# Row 27
$ws.Range("H27").Value = 1622.2941
$ws.Range("I27").Value = 369.33334
$ws.Range("K27").Value = 369.33334
$ws.Range("M27").Value = -262.33334
# Row 55
$ws.Range("H55").Value = 4996.8
$ws.Range("J55").Value = 4777.6665
$ws.Range("L55").Value = 4777.6665
$ws.Range("N55").Value = -5123.6665
# Row 132
$ws.Range("H132").Value = 8488.125
$ws.Range("J132").Value = 12315.857
$ws.Range("L132").Value = 36947.571
$ws.Range("N132").Value = -42007.571
# Row 136
$ws.Range("H136").Value = 2966.2632
$ws.Range("I136").Value = 2131
$ws.Range("J136").Value = 6098.5
$ws.Range("K136").Value = 6393
$ws.Range("L136").Value = 18295.5
$ws.Range("M136").Value = -3843
$ws.Range("N136").Value = -23395.5

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1216.3334
$ws.Range("I100").Value = 1299.6
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 2599.2
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -2058.2
$ws.Range("N100").Value = -2682
# Row 107
$ws.Range("H107").Value = 1506.9
$ws.Range("I107").Value = 1708.75
$ws.Range("K107").Value = 5126.25
$ws.Range("M107").Value = -3206.25
# Row 122
$ws.Range("H122").Value = 2649.875
$ws.Range("I122").Value = 2618.0454
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7854.1362
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5404.1362
$ws.Range("N122").Value = -13900
# Row 132
$ws.Range("H132").Value = 2042.5454
$ws.Range("I132").Value = 2042.5454
$ws.Range("K132").Value = 6127.6362
$ws.Range("M132").Value = -3597.6362
